$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header "time_taken" in F1, copying the header style from E1 first
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Fill F2:F41 with the recorded time_taken values (plain text, default style)
$timeTaken = @(
    "2021-10-05 13:38:37.254814",
    "2021-10-05 13:38:37.254825",
    "2021-10-05 13:38:37.254829",
    "2021-10-05 13:38:37.254832",
    "2021-10-05 13:38:37.254836",
    "2021-10-05 13:38:37.254839",
    "2021-10-05 13:38:37.254842",
    "2021-10-05 13:38:37.254845",
    "2021-10-05 13:38:37.254848",
    "2021-10-05 13:38:37.254851",
    "2021-10-05 13:38:37.254854",
    "2021-10-05 13:38:37.254857",
    "2021-10-05 13:38:37.254860",
    "2021-10-05 13:38:37.254863",
    "2021-10-05 13:38:37.254866",
    "2021-10-05 13:38:37.254869",
    "2021-10-05 13:38:37.254872",
    "2021-10-05 13:38:37.254875",
    "2021-10-05 13:38:37.254878",
    "2021-10-05 13:38:37.254881",
    "2021-10-05 13:38:37.254884",
    "2021-10-05 13:38:37.254887",
    "2021-10-05 13:38:37.254890",
    "2021-10-05 13:38:37.254893",
    "2021-10-05 13:38:37.254896",
    "2021-10-05 13:38:37.254899",
    "2021-10-05 13:38:37.254903",
    "2021-10-05 13:38:37.254905",
    "2021-10-05 13:38:37.254908",
    "2021-10-05 13:38:37.254911",
    "2021-10-05 13:38:37.254914",
    "2021-10-05 13:38:37.254917",
    "2021-10-05 13:38:37.254920",
    "2021-10-05 13:38:37.254923",
    "2021-10-05 13:38:37.254926",
    "2021-10-05 13:38:37.254929",
    "2021-10-05 13:38:37.254932",
    "2021-10-05 13:38:37.254935",
    "2021-10-05 13:38:37.254938",
    "2021-10-05 13:38:37.254941"
)

for ($i = 0; $i -lt $timeTaken.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

Write-Host "done"
